# "Refactored and added Tables"
#
# 1. Rename sheet1 "Add Custome" -> "Add Customer"
# 2. Rework the "Add Customer" sheet so column G builds a single bulk
#    "insert into customer(...)values" statement (header in G2) followed by
#    one row-tuple per data row (G3:G29), instead of one full INSERT per row.
# 3. Add two more (currently empty) sheets: "Add Invoice" and "Add Product".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Add Customer"

# --- Header row: build "insert into customer(...)values" in G2 ---
$ws.Range("G2").Formula = '="insert into customer("&B2&","&C2&","&D2&","&E2&","&F2&")values"'

# --- Serial numbers for the two existing data rows ---
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2

# --- Existing rows 3 & 4: switch from full INSERT statements to
#     comma-terminated value tuples ---
$ws.Range("G3").Formula = '="(''"&B3&"'',''"&C3&"'',''"&D3&"'',''"&E3&"'',''"&F3&"''),"'
$ws.Range("G4").Formula = '="(''"&B4&"'',''"&C4&"'',''"&D4&"'',''"&E4&"'',''"&F4&"''),"'

# --- Drop the leftover literal text that used to sit in B9 ---
$ws.Range("B9").ClearContents() | Out-Null

# --- F5 only carries the (empty) Hyperlink style, like F3/F4 ---
$ws.Range("F5").Style = "Hyperlink"

# --- Fill G5:G29 with the shared row-tuple formula (blank placeholder rows
#     ready for more customers to be typed into A:F) ---
$ws.Range("G5:G29").Formula = '="(''"&B5&"'',''"&C5&"'',''"&D5&"'',''"&E5&"'',''"&F5&"''),"'

# --- Column sizing to fit the new layout (S/N, Firstname..Email, generated
#     SQL column) ---
$ws.Columns.Item(1).ColumnWidth = 4
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
$ws.Columns.Item(4).ColumnWidth = 20.333333333333332
$ws.Columns.Item(7).ColumnWidth = 97.66666666666667

# --- Restore a normal view (no frozen/offset top-left cell) with the
#     selection parked at B6, the first blank entry row ---
$ws.Range("B6").Select() | Out-Null

# --- Add the two new (blank) sheets after "Add Customer" ---
$wsInvoice = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsInvoice.Name = "Add Invoice"
$wsInvoice.Range("G26").Select() | Out-Null

$wsProduct = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsProduct.Name = "Add Product"
$wsProduct.Range("F28").Select() | Out-Null

# --- Leave "Add Customer" as the active/selected sheet ---
$ws.Activate() | Out-Null
